$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Glavni tok")
$ws1.Range("C3").Value = "2. Skenira boarding pass"
